# Updated cryptos list on Tue Jun 11 11:30:15 UTC 2024 with GitHub Actions
#
# Prices/volumes are stored as plain text in this sheet (no number format),
# so a handful of the new "Price" values look like valid numbers to Excel's
# normal Range.Value auto-detection (e.g. "605.38", "1.00"). For those we
# prefix the literal with a leading apostrophe, exactly as a user would type
# it in the Excel UI, to force a text entry instead of a numeric one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.027.74"
$ws.Range("E2").Value = "  -3.40%  "

$ws.Range("D3").Value = "3.536.51"
$ws.Range("E3").Value = "  -3.76%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'605.38"
$ws.Range("E5").Value = "  -5.62%  "

$ws.Range("D6").Value = "'154.50"
$ws.Range("E6").Value = "  -3.20%  "

$ws.Range("D7").Value = "3.534.41"
$ws.Range("E7").Value = "  -3.76%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("D11").Value = "'6.83"
$ws.Range("E11").Value = "  -3.38%  "

$ws.Range("D12").Value = "'0.430"
$ws.Range("E12").Value = "  -3.71%  "

$ws.Range("E13").Value = "  -4.35%  "

$ws.Range("D14").Value = "4.132.89"
$ws.Range("E14").Value = "  -3.79%  "

$ws.Range("D15").Value = "'31.99"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "3.545.80"
$ws.Range("E16").Value = "  -3.04%  "

$ws.Range("D17").Value = "67.054.35"
$ws.Range("E17").Value = "  -3.37%  "

$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "'6.37"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("D20").Value = "'15.47"
$ws.Range("E20").Value = "  -3.17%  "

$ws.Range("D21").Value = "'451.92"

$ws.Range("D22").Value = "'9.37"
$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D24").Value = "'79.02"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").Value = "3.675.02"
$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").Value = "'10.24"
$ws.Range("E28").Value = "  -6.01%  "

$ws.Range("E29").Value = "  -8.15%  "

$ws.Range("D30").Value = "'1.68"
$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("E31").Value = "  -2.89%  "

$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").Value = "'25.95"
$ws.Range("E33").Value = "  -3.42%  "

$ws.Range("E34").Value = "  -5.20%  "

$ws.Range("E35").Value = "  -3.91%  "

$ws.Range("E36").Value = "  -4.89%  "

$ws.Range("D37").Value = "3.528.81"
$ws.Range("E37").Value = "  -3.79%  "

$ws.Range("D38").Value = "'8.09"
$ws.Range("E38").Value = "  -4.10%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "'176.45"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("D43").Value = "'5.60"
$ws.Range("E43").Value = "  -4.86%  "

$ws.Range("D44").Value = "'0.0876"
$ws.Range("E44").Value = "  -2.49%  "

$ws.Range("D45").Value = "'0.893"
$ws.Range("E45").Value = "  -3.47%  "

$ws.Range("D46").Value = "'45.80"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").Value = "'28.47"
$ws.Range("E47").Value = "  +4.15%  "

$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("E49").Value = "  -1.42%  "

# Rows 50 and 51 swapped coins: Cosmos <-> SuiNetwork, with new price/volume data.
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'1.04"
$ws.Range("E50").Value = "  -3.34%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "  -2.32%  "
